# Auto-generated edit script applying the Masamune_Profits diff
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H8").Value = 78.181816
$ws.Range("I8").Value = 66
$ws.Range("J8").Value = 200
$ws.Range("K8").Value = 198
$ws.Range("L8").Value = 600
$ws.Range("M8").Value = -59
$ws.Range("N8").Value = -878
$ws.Range("H98").Value = 27316.434
$ws.Range("I98").Value = 1006.3158
$ws.Range("J98").Value = 72761.17999999999
$ws.Range("K98").Value = 1006.3158
$ws.Range("L98").Value = 72761.17999999999
$ws.Range("M98").Value = 491.6842
$ws.Range("N98").Value = -75757.17999999999
$ws.Range("H100").Value = 2748.5
$ws.Range("I100").Value = 2995
$ws.Range("J100").Value = 2649.9
$ws.Range("K100").Value = 2995
$ws.Range("L100").Value = 2649.9
$ws.Range("M100").Value = -2454
$ws.Range("N100").Value = -3731.9
$ws.Range("H113").Value = 2913.8572
$ws.Range("I113").Value = 2899.25
$ws.Range("J113").Value = 2933.3333
$ws.Range("K113").Value = 2899.25
$ws.Range("L113").Value = 2933.3333
$ws.Range("M113").Value = 354.75
$ws.Range("N113").Value = -9441.3333
$ws.Range("H122").Value = 27316.434
$ws.Range("I122").Value = 1006.3158
$ws.Range("J122").Value = 72761.17999999999
$ws.Range("K122").Value = 3018.9474
$ws.Range("L122").Value = 218283.54
$ws.Range("M122").Value = -568.9474
$ws.Range("N122").Value = -223183.54
$ws.Range("H138").Value = 2919.6667
$ws.Range("I138").Value = 2508.4348
$ws.Range("J138").Value = 3197.853
$ws.Range("K138").Value = 7525.3044
$ws.Range("L138").Value = 9593.559000000001
$ws.Range("M138").Value = -2385.3044
$ws.Range("N138").Value = -19873.559

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 1714
$ws.Range("I122").Value = 1349.5
$ws.Range("J122").Value = 1859.8
$ws.Range("K122").Value = 4048.5
$ws.Range("L122").Value = 5579.4
$ws.Range("M122").Value = -1598.5
$ws.Range("N122").Value = -10479.4
$ws.Range("H125").Value = 50715
$ws.Range("J125").Value = 50715
$ws.Range("L125").Value = 50715
$ws.Range("N125").Value = -60555
$ws.Range("H132").Value = 13891406
$ws.Range("I132").Value = 16668871
$ws.Range("J132").Value = 4082.3333
$ws.Range("K132").Value = 50006613
$ws.Range("L132").Value = 12246.9999
$ws.Range("M132").Value = -50004083
$ws.Range("N132").Value = -17306.9999
$ws.Range("H138").Value = 51210
$ws.Range("J138").Value = 51210
$ws.Range("L138").Value = 51210
$ws.Range("N138").Value = -61490

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H130").Value = 49945
$ws.Range("J130").Value = 49945
$ws.Range("L130").Value = 49945
$ws.Range("N130").Value = -59985
$ws.Range("H134").Value = 4573.5
$ws.Range("I134").Value = 4664
$ws.Range("J134").Value = 4519.2
$ws.Range("K134").Value = 13992
$ws.Range("L134").Value = 13557.6
$ws.Range("M134").Value = -11457
$ws.Range("N134").Value = -18627.6

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 8011.926
$ws.Range("I31").Value = 3723.3333
$ws.Range("J31").Value = 10156.223
$ws.Range("K31").Value = 3723.3333
$ws.Range("L31").Value = 10156.223
$ws.Range("M31").Value = -3428.3333
$ws.Range("N31").Value = -10746.223
$ws.Range("H34").Value = 8011.926
$ws.Range("I34").Value = 3723.3333
$ws.Range("J34").Value = 10156.223
$ws.Range("K34").Value = 3723.3333
$ws.Range("L34").Value = 10156.223
$ws.Range("M34").Value = -3521.3333
$ws.Range("N34").Value = -10560.223
$ws.Range("H122").Value = 83436010
$ws.Range("I122").Value = 111245416
$ws.Range("J122").Value = 7804.6665
$ws.Range("K122").Value = 333736248
$ws.Range("L122").Value = 23413.9995
$ws.Range("M122").Value = -333733798
$ws.Range("N122").Value = -28313.9995

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 2812.0334
$ws.Range("J131").Value = 1191.509
$ws.Range("L131").Value = 3574.527
$ws.Range("N131").Value = -13654.527

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 935.6667
$ws.Range("I122").Value = 903.5
$ws.Range("J122").Value = 1000
$ws.Range("K122").Value = 2710.5
$ws.Range("L122").Value = 3000
$ws.Range("M122").Value = -260.5
$ws.Range("N122").Value = -7900

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 66668908
$ws.Range("I7").Value = 111112790
$ws.Range("J7").Value = 3084.1667
$ws.Range("K7").Value = 111112790
$ws.Range("L7").Value = 3084.1667
$ws.Range("M7").Value = -111112678
$ws.Range("N7").Value = -3308.1667
$ws.Range("H16").Value = 1064.5807
$ws.Range("I16").Value = 1064.5807
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 1064.5807
$ws.Range("L16").Value = 0
$ws.Range("M16").Value = -894.5807
$ws.Range("N16").ClearContents()
$ws.Range("H122").Value = 64412.75
$ws.Range("I122").Value = 73307.42999999999
$ws.Range("J122").Value = 2150
$ws.Range("K122").Value = 219922.29
$ws.Range("L122").Value = 6450
$ws.Range("M122").Value = -217472.29
$ws.Range("N122").Value = -11350
$ws.Range("H126").Value = 66668908
$ws.Range("I126").Value = 111112790
$ws.Range("J126").Value = 3084.1667
$ws.Range("K126").Value = 333338370
$ws.Range("L126").Value = 9252.500100000001
$ws.Range("M126").Value = -333335900
$ws.Range("N126").Value = -14192.5001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H6").Value = 1390248.9
$ws.Range("I6").Value = 11111111
$ws.Range("J6").Value = 1554.2858
$ws.Range("K6").Value = 11111111
$ws.Range("L6").Value = 11111111
$ws.Range("M6").Value = -11110996
$ws.Range("N6").Value = -1784.2858
$ws.Range("H13").Value = 1160
$ws.Range("I13").Value = 1000
$ws.Range("K13").Value = 1000
$ws.Range("M13").Value = -860
$ws.Range("H64").Value = 11000
$ws.Range("J64").Value = 11000
$ws.Range("L64").Value = 11000
$ws.Range("N64").Value = -11496
$ws.Range("H67").Value = 11000
$ws.Range("J67").Value = 11000
$ws.Range("L67").Value = 11000
$ws.Range("N67").Value = -12716
$ws.Range("H122").Value = 31093372
$ws.Range("I122").Value = 48052812
$ws.Range("J122").Value = 1064.8334
$ws.Range("K122").Value = 144158436
$ws.Range("L122").Value = 3194.5002
$ws.Range("M122").Value = -144155986
$ws.Range("N122").Value = -8094.5002
$ws.Range("H126").Value = 1961631.9
$ws.Range("I126").Value = 5882855
$ws.Range("J126").Value = 1020.2
$ws.Range("K126").Value = 17648565
$ws.Range("L126").Value = 3060.6
$ws.Range("M126").Value = -17646095
$ws.Range("N126").Value = -8000.6
$ws.Range("H132").Value = 1550.4857
$ws.Range("I132").Value = 1234.6129
$ws.Range("J132").Value = 3998.5
$ws.Range("K132").Value = 3703.8387
$ws.Range("L132").Value = 11995.5
$ws.Range("M132").Value = -1173.8387
$ws.Range("N132").Value = -17055.5
$ws.Range("H136").Value = 18113.115
$ws.Range("I136").Value = 30511.766
$ws.Range("J136").Value = 2500
$ws.Range("K136").Value = 91535.298
$ws.Range("L136").Value = 7500
$ws.Range("M136").Value = -88985.298
$ws.Range("N136").Value = -12600
